# This script updates column C (D_lambda [ppm]) values in the worksheet
# for rows 2-65, reflecting a re-run of the transit-depth simulation with
# different temperature parameters ("Simulações com diferentes temperaturas").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 698.342243854011
$ws.Range("C3").Value = 593.6769593554159
$ws.Range("C4").Value = 593.0411156351845
$ws.Range("C5").Value = 666.4419341200212
$ws.Range("C6").Value = 665.1522464695691
$ws.Range("C7").Value = 680.6202059079781
$ws.Range("C8").Value = 736.3712711970783
$ws.Range("C9").Value = 702.319849793187
$ws.Range("C10").Value = 688.2115340245409
$ws.Range("C11").Value = 718.6013598591279
$ws.Range("C12").Value = 721.1324989108148
$ws.Range("C13").Value = 708.8691411460957
$ws.Range("C14").Value = 714.8109949277615
$ws.Range("C15").Value = 718.6004537045321
$ws.Range("C16").Value = 737.2757459551682
$ws.Range("C17").Value = 750.6694883916421
$ws.Range("C18").Value = 768.3922218387762
$ws.Range("C19").Value = 763.1498285415316
$ws.Range("C20").Value = 764.5814621501934
$ws.Range("C21").Value = 770.5334982478496
$ws.Range("C22").Value = 779.8273020018787
$ws.Range("C23").Value = 789.9509190403276
$ws.Range("C24").Value = 799.6833002618109
$ws.Range("C25").Value = 805.5469223620104
$ws.Range("C26").Value = 811.8072862606241
$ws.Range("C27").Value = 817.6684168161863
$ws.Range("C28").Value = 826.3841619271517
$ws.Range("C29").Value = 831.5088430658202
$ws.Range("C30").Value = 834.8227068875058
$ws.Range("C31").Value = 841.1974692998924
$ws.Range("C32").Value = 844.8471400407387
$ws.Range("C33").Value = 848.3788786642909
$ws.Range("C34").Value = 851.4976800779727
$ws.Range("C35").Value = 853.9148473436509
$ws.Range("C36").Value = 856.0605157041934
$ws.Range("C37").Value = 858.6283344400591
$ws.Range("C38").Value = 859.9942350322465
$ws.Range("C39").Value = 862.0399480187002
$ws.Range("C40").Value = 865.0750422755093
$ws.Range("C41").Value = 867.4083611793782
$ws.Range("C42").Value = 867.6301866678004
$ws.Range("C43").Value = 869.2860803497515
$ws.Range("C44").Value = 870.4101049136215
$ws.Range("C45").Value = 871.3663389674498
$ws.Range("C46").Value = 871.5631265351353
$ws.Range("C47").Value = 872.2474665332492
$ws.Range("C48").Value = 873.3546947625737
$ws.Range("C49").Value = 873.8045907060732
$ws.Range("C50").Value = 875.0707375694278
$ws.Range("C51").Value = 876.8854149621142
$ws.Range("C52").Value = 878.911688527273
$ws.Range("C53").Value = 882.6043870351352
$ws.Range("C54").Value = 886.1543106184478
$ws.Range("C55").Value = 891.8155747982271
$ws.Range("C56").Value = 897.2423174790079
$ws.Range("C57").Value = 899.5856096751797
$ws.Range("C58").Value = 900.4707647216615
$ws.Range("C59").Value = 901.5686295065217
$ws.Range("C60").Value = 904.656407209603
$ws.Range("C61").Value = 911.2034788496537
$ws.Range("C62").Value = 912.1673648400108
$ws.Range("C63").Value = 913.1174061474123
$ws.Range("C64").Value = 913.7790811927226
$ws.Range("C65").Value = 914.6378137127042
